$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "BRIDLEFIELD LANE MARKHAM ON L6C2P3 61"
$ws.Range("A1").Select()
